# Auto-generated edit script: updates cryptos list values per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking Price (D) values are written as text, matching the
# original inlineStr cell type (not auto-converted to numbers by Excel).
$dRefs = @("D2", "D3", "D5", "D6", "D8", "D11", "D12", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D24", "D27", "D29", "D32", "D33", "D34", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D48", "D50", "D51")
foreach ($r in $dRefs) { $ws.Range($r).NumberFormat = "@" }

$ws.Range("D2").Value = "62.916.55"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "3.167.92"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "589.10"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").Value = "138.50"
$ws.Range("E6").Value = "  -3.15%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.162.48"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").Value = "5.34"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").Value = "0.458"
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("E13").Value = "  -3.90%  "
$ws.Range("D14").Value = "33.98"
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("D15").Value = "3.685.61"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "3.164.87"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "62.894.53"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "6.66"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").Value = "471.72"
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").Value = "13.92"
$ws.Range("E21").Value = "  -5.14%  "
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").Value = "83.57"
$ws.Range("E24").Value = "  -3.68%  "
$ws.Range("E25").Value = "  -3.38%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "2.70"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D29").Value = "7.94"
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "26.76"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").Value = "0.108"
$ws.Range("E33").Value = "  -2.94%  "
$ws.Range("D34").Value = "2.51"
$ws.Range("E34").Value = "  -4.80%  "
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("D37").Value = "5.77"
$ws.Range("E37").Value = "  -3.71%  "
$ws.Range("D38").Value = "0.0₃0704"
$ws.Range("E38").Value = "  -5.40%  "
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("D40").Value = "415.63"
$ws.Range("E40").Value = "  -5.04%  "
$ws.Range("D41").Value = "2.955.77"
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.74"
$ws.Range("E42").Value = "  -7.67%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "8.29"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "0.110"
$ws.Range("E44").Value = "  -8.12%  "
$ws.Range("D45").Value = "0.262"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("D48").Value = "25.47"
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "2.24"
$ws.Range("E50").Value = "  -5.55%  "
$ws.Range("D51").Value = "119.41"
$ws.Range("E51").Value = "  -1.69%  "

Write-Host "Updated $($dRefs.Count) price cells and related fields."
